{"js": "// Word JS API (Office.js) edit script.\n//\n// Commit message: \"clockwise corrected to anti-clockwise in all documents\"\n//\n// The testing-report table has a row reading:\n//   \"5.2 All counters are redistributed to other pits in clockwise direction.\"\n// which must be corrected to:\n//   \"5.2 All counters are redistributed to other pits in anti-clockwise direction.\"\n//\n// Find the stand-alone word \"clockwise\" and turn it into \"anti-clockwise\".\n// Word.InsertLocation.replace swaps the text of the matched range in place,\n// so the surrounding run formatting (the green highlight / en-GB language\n// tag) is preserved untouched.\n\nconst body = context.document.body;\n\nconst results = body.search(\"clockwise\", { matchCase: true, matchWholeWord: true });\nresults.load(\"text\");\nawait context.sync();\n\nfor (const result of results.items) {\n  // \"clockwise\" is also matched as a whole word inside \"anti-clockwise\"\n  // (hyphen counts as a word boundary), so guard against double-applying\n  // the fix by checking the enclosing paragraph text first.\n  const para = result.paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n\n  if (para.text.indexOf(\"anti-clockwise\") === -1) {\n    result.insertText(\"anti-clockwise\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Commit message: \"clockwise corrected to anti-clockwise in all documents\"\n#\n# The testing-report table has a row reading:\n#   \"5.2 All counters are redistributed to other pits in clockwise direction.\"\n# which must be corrected to:\n#   \"5.2 All counters are redistributed to other pits in anti-clockwise direction.\"\n\n$d = $word.ActiveDocument\n\n# Idempotency guard: only touch the document if the fix hasn't been applied yet.\nif (-not $d.Content.Text.Contains(\"anti-clockwise\")) {\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $find.Text = \"clockwise\"\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.Replacement.Text = \"anti-clockwise\"\n\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(wdFindContinue=1), Format,\n    # ReplaceWith, Replace(wdReplaceAll=2)\n    $find.Execute(\"clockwise\", $true, $true, $false, $false, $false, $true, 1, $false, \"anti-clockwise\", 2) | Out-Null\n}\n"}
